# Fruta / hortaliza, semanal
# Insert a new weekly record at row 27 (Feria Lagunitas de Puerto Montt - Arveja Verde),
# shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 27; all rows from 27 downward shift to 28..91,
# carrying their existing values/styles with them.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record.
$ws.Range("A27").Value = 4
$ws.Range("B27").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C27").Value = 'Los Lagos'
$ws.Range("D27").Value = 44526
$ws.Range("E27").Value = 10
$ws.Range("F27").Value = 100112022
$ws.Range("G27").Value = 'Arveja Verde'
$ws.Range("H27").Value = 'Sin especificar'
$ws.Range("I27").Value = 'Primera'
$ws.Range("J27").Value = 80
$ws.Range("K27").Value = 17000
$ws.Range("L27").Value = 17000
$ws.Range("M27").Value = 17000
$ws.Range("N27").Value = '$/saco 25 kilos'
$ws.Range("O27").Value = 'Región del Maule'
$ws.Range("P27").Value = 680
$ws.Range("Q27").Value = 25
$ws.Range("R27").Value = 'Hortaliza'
